$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "dzeca"
$ws.Range("B5").Value = "dzeca"
$ws.Range("C5").Value = "dzeca@gmail"
$ws.Range("D5").Value = ";sd;fklad;klf"

# "20" looks numeric, but the source row keeps it as text (shared string),
# same as the other numeric-looking entries already in the sheet. Force the
# text type for the write, then restore the default (unstyled) look so the
# cell doesn't end up carrying an explicit style like the rest of the row.
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "20"
$ws.Range("E5").Style = "Normal"
